$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "nombre" (B) column with "SMG" for rows 2-4
$ws.Range("B2").Value = "SMG"
$ws.Range("B3").Value = "SMG"
$ws.Range("B4").Value = "SMG"

# Fill in the "parrafo" (D) column with text for rows 2-4
$ws.Range("D2").Value = "Hola somos talalalsaldafasfalsf"
$ws.Range("D3").Value = "asfsagfasdgadsgasdgadsgasdgadg"
$ws.Range("D4").Value = "agadsgadgagadgadgasfrwqfSCs"

# Update the selection to D4
$ws.Range("D4").Select()
